# MP1: Better dB processing
# Applies the edit described by the commit: extends the Scales table with a
# new data point, adds three new "3750/1250 cent" rows + derived chart
# series, and repositions/resizes the three charts to make room.

$wb = $excel.ActiveWorkbook
$wsInstr = $wb.Worksheets.Item("Instr")
$ws = $wb.Worksheets.Item("Scales")

# ---------------------------------------------------------------------------
# 1. Scales!row1/row5 header rows: existing P gets a new value, old P value
#    shifts out to the new Q column.
# ---------------------------------------------------------------------------
$ws.Range("P1").Copy($ws.Range("Q1"))
$ws.Range("P1").Value = 1

$ws.Range("P5").Copy($ws.Range("Q5"))
$ws.Range("P5").ClearContents()

# ---------------------------------------------------------------------------
# 2. Scales!row2/3/4: the old "-inf" marker cell (P) moves to Q, and P gets a
#    real new numeric data point.
# ---------------------------------------------------------------------------
$ws.Range("P2").Copy($ws.Range("Q2"))
$ws.Range("O2").Copy()
$ws.Range("P2").PasteSpecial(-4122)
$ws.Range("P2").Value = -113.8

$ws.Range("P3").Copy($ws.Range("Q3"))
$ws.Range("P3").ClearFormats()
$ws.Range("P3").Value = -125.8

$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("P4").ClearFormats()
$ws.Range("P4").Value = -42

# ---------------------------------------------------------------------------
# 3. Scales!row6-11: simple rename of the trailing literal/formula cell from
#    column P to column Q (no new data point here, the table just grew).
# ---------------------------------------------------------------------------
$ws.Range("P6").Copy($ws.Range("Q6"))
$ws.Range("P6").Clear()

$ws.Range("P7").Copy($ws.Range("Q7"))
$ws.Range("P7").Clear()

$ws.Range("P8").Copy($ws.Range("Q8"))
$ws.Range("P8").Clear()

$ws.Range("P9").Copy($ws.Range("Q9"))
$ws.Range("P9").Clear()

$ws.Range("P10").Copy($ws.Range("Q10"))
$ws.Range("P10").Clear()

$ws.Range("P11").Copy($ws.Range("Q11"))
$ws.Range("Q11").Formula = "=POWER(10, Q10/20)"
$ws.Range("P11").Clear()

# ---------------------------------------------------------------------------
# 4. New rows 12-14: three new "cent deviation" curves driven off row 1, plus
#    helper rows 15 and 20.
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = 0.007
$ws.Range("B12").Value = 3750
$ws.Range("C12").Value = 0
$ws.Range("D12").Formula = "=(POWER(D`$1/127, `$A12)-1)*`$B12"
$ws.Range("E12:P12").Formula = "=(POWER(E`$1/127, `$A12)-1)*`$B12"

$ws.Range("A13").Value = 0.007
$ws.Range("B13").Formula = "=B12/3"
$ws.Range("C13").Value = 0
$ws.Range("D13").Formula = "=(POWER(D`$1/127, `$A13)-1)*`$B13"
$ws.Range("E13:P13").Formula = "=(POWER(E`$1/127, `$A13)-1)*`$B13"

$ws.Range("A14").Value = 0.007
$ws.Range("B14").Value = 3750
$ws.Range("C14").Value = 12
$ws.Range("D14").Formula = "=(POWER(D`$1/127, `$A14)-1)*`$B14+12"
$ws.Range("E14:P14").Formula = "=(POWER(E`$1/127, `$A14)-1)*`$B14+12"

$ws.Range("A15").Formula = "=1/A14"

$ws.Range("P20").Formula = "=(POWER(111.75/127, 0.007) - 1) *1250*2"

# ---------------------------------------------------------------------------
# 5. Chart 5 ("Piano/Violin/Flute" vs full CC scale, Scales!C5:*5 based):
#    every series range grows from column P to column Q.
# ---------------------------------------------------------------------------
$chart5 = $ws.ChartObjects().Item(1).Chart
$chart5.SeriesCollection().Item(1).Formula = "=SERIES(Scales!`$B`$6,Scales!`$C`$5:`$Q`$5,Scales!`$C`$6:`$Q`$6,1)"
$chart5.SeriesCollection().Item(2).Formula = "=SERIES(Scales!`$B`$7,Scales!`$C`$5:`$Q`$5,Scales!`$C`$7:`$Q`$7,2)"
$chart5.SeriesCollection().Item(3).Formula = "=SERIES(Scales!`$B`$8,Scales!`$C`$5:`$Q`$5,Scales!`$C`$8:`$Q`$8,3)"

# ---------------------------------------------------------------------------
# 6. Chart 6 ("CC to db", Scales!C1:*1 based): ranges grow by one column, and
#    three brand-new series (rows 12/13/14 - the "3750/1250 cent" curves)
#    get appended.
# ---------------------------------------------------------------------------
$chart6 = $ws.ChartObjects().Item(2).Chart
$chart6.SeriesCollection().Item(1).Formula = "=SERIES(Scales!`$B`$2,Scales!`$C`$1:`$P`$1,Scales!`$C`$2:`$P`$2,1)"
$chart6.SeriesCollection().Item(2).Formula = "=SERIES(Scales!`$B`$3,Scales!`$C`$1:`$P`$1,Scales!`$C`$3:`$P`$3,2)"
$chart6.SeriesCollection().Item(3).Formula = "=SERIES(Scales!`$B`$4,Scales!`$C`$1:`$P`$1,Scales!`$C`$4:`$P`$4,3)"

$s4 = $chart6.SeriesCollection().NewSeries()
$s4.Formula = "=SERIES(Scales!`$B`$12,Scales!`$C`$1:`$P`$1,Scales!`$C`$12:`$P`$12,4)"

$s5 = $chart6.SeriesCollection().NewSeries()
$s5.Formula = "=SERIES(Scales!`$B`$13,,Scales!`$C`$13:`$P`$13,5)"

$s6 = $chart6.SeriesCollection().NewSeries()
$s6.Formula = "=SERIES(Scales!`$B`$14,,Scales!`$C`$14:`$P`$14,6)"

# ---------------------------------------------------------------------------
# 7. Chart 7 ("Amplitude", Scales!D10:*10 based): ranges grow by one column.
# ---------------------------------------------------------------------------
$chart7 = $ws.ChartObjects().Item(3).Chart
$chart7.SeriesCollection().Item(1).Formula = "=SERIES(Scales!`$C`$11,Scales!`$D`$10:`$Q`$10,Scales!`$D`$11:`$Q`$11,1)"

# ---------------------------------------------------------------------------
# 8. Re-anchor/resize the three charts to their new drawing positions. All
#    columns/rows on this sheet use the default size, so the EMU anchors
#    from the target drawing XML can be reproduced exactly in points.
# ---------------------------------------------------------------------------
function ColToPoints($sheet, $colIndex0, $offEmu) {
    $cum = 0.0
    for ($c = 1; $c -le $colIndex0; $c++) {
        $cum += $sheet.Columns.Item($c).Width
    }
    return $cum + ($offEmu / 12700.0)
}
function RowToPoints($sheet, $rowIndex0, $offEmu) {
    $cum = 0.0
    for ($r = 1; $r -le $rowIndex0; $r++) {
        $cum += $sheet.Rows.Item($r).Height
    }
    return $cum + ($offEmu / 12700.0)
}

$co5 = $ws.ChartObjects().Item(1)
$fromX = ColToPoints $ws 18 18097
$fromY = RowToPoints $ws 0 38099
$toX = ColToPoints $ws 30 66675
$toY = RowToPoints $ws 40 133349
$co5.Left = $fromX
$co5.Top = $fromY
$co5.Width = $toX - $fromX
$co5.Height = $toY - $fromY

$co6 = $ws.ChartObjects().Item(2)
$fromX = ColToPoints $ws 2 35241
$fromY = RowToPoints $ws 14 139065
$toX = ColToPoints $ws 14 300990
$toY = RowToPoints $ws 40 38100
$co6.Left = $fromX
$co6.Top = $fromY
$co6.Width = $toX - $fromX
$co6.Height = $toY - $fromY

$co7 = $ws.ChartObjects().Item(3)
$fromX = ColToPoints $ws 30 98107
$fromY = RowToPoints $ws 0 39052
$toX = ColToPoints $ws 37 402907
$toY = RowToPoints $ws 16 67627
$co7.Left = $fromX
$co7.Top = $fromY
$co7.Width = $toX - $fromX
$co7.Height = $toY - $fromY

# ---------------------------------------------------------------------------
# 9. Active sheet / selection bookkeeping: "Instr" becomes the active tab.
# ---------------------------------------------------------------------------
$ws.Range("P20").Select()
$wsInstr.Range("F16").Select()
$wsInstr.Activate()
